$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Liberar la pantalla (fila 7): quitar el usuario y marcar como "Libre"
$ws.Range("B7").ClearContents()
$ws.Range("C7").Value = "Libre"
